$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.385.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.883.72'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08034'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3126'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08333'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.891.66'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7211'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.48%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.249'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.92'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.339'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008449'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.399.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('B19').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C19').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.152.90'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '241.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.863'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1589'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.98'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.048'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.56'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.419'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.199'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.27%  '
$ws.Range('E33').Value = '  +2.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.951'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7503'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.702'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.291.89'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.27%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01889'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.748'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.602'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9202'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '111.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '74.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.69%  '
$ws.Range('E46').Value = '  +5.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.039.68'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.52%  '
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5219'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.510'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4406'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.21%  '
